$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 151.7260716666667
$ws.Range("H2").Value = 455.178215
$ws.Range("I2").Value = 0.2700739458961593
$ws.Range("J2").Value = 0.2783366498663096
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.275247666666667
$ws.Range("N2").Value = 24.825743
$ws.Range("O2").Value = 0.1035475654300768
$ws.Range("P2").Value = 0.1152196269592993
$ws.Range("Q2").Value = 1255.570820532083
$ws.Range("R2").Value = 11300.13738478874
$ws.Range("S2").Value = 0.02796549958364157
$ws.Range("T2").Value = 0.03206984496669728

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 151.7260716666667
$ws.Range("H3").Value = 455.178215
$ws.Range("I3").Value = 0.2700739458961593
$ws.Range("J3").Value = 0.2783366498663096
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.354572
$ws.Range("N3").Value = 142.063716
$ws.Range("O3").Value = 0.5925442766305059
$ws.Range("P3").Value = 0.6593368972671567
$ws.Range("Q3").Value = 7184.923185016327
$ws.Range("R3").Value = 64664.30866514694
$ws.Range("S3").Value = 0.1600307709077861
$ws.Range("T3").Value = 0.1835176231185875

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 151.7260716666667
$ws.Range("H4").Value = 455.178215
$ws.Range("I4").Value = 0.2700739458961593
$ws.Range("J4").Value = 0.2783366498663096
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 24.2875365
$ws.Range("N4").Value = 48.575073
$ws.Range("O4").Value = 0.3039081579394173
$ws.Range("P4").Value = 0.225443475773544
$ws.Range("Q4").Value = 3685.052503605782
$ws.Range("R4").Value = 22110.31502163469
$ws.Range("S4").Value = 0.08207767540473161
$ws.Range("T4").Value = 0.06274918178102476

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 82.248871
$ws.Range("H5").Value = 246.746613
$ws.Range("I5").Value = 0.146403824289839
$ws.Range("J5").Value = 0.150882936320401
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.275247666666667
$ws.Range("N5").Value = 24.825743
$ws.Range("O5").Value = 0.1035475654300768
$ws.Range("P5").Value = 0.1152196269592993
$ws.Range("Q5").Value = 680.6297778287177
$ws.Range("R5").Value = 6125.668000458459
$ws.Range("S5").Value = 0.01515975957486557
$ws.Range("T5").Value = 0.01738467563736031

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 82.248871
$ws.Range("H6").Value = 246.746613
$ws.Range("I6").Value = 0.146403824289839
$ws.Range("J6").Value = 0.150882936320401
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.354572
$ws.Range("N6").Value = 142.063716
$ws.Range("O6").Value = 0.5925442766305059
$ws.Range("P6").Value = 0.6593368972671567
$ws.Range("Q6").Value = 3894.860083688211
$ws.Range("R6").Value = 35053.74075319391
$ws.Range("S6").Value = 0.08675074815976233
$ws.Range("T6").Value = 0.09948268708405117

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 82.248871
$ws.Range("H7").Value = 246.746613
$ws.Range("I7").Value = 0.146403824289839
$ws.Range("J7").Value = 0.150882936320401
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 24.2875365
$ws.Range("N7").Value = 48.575073
$ws.Range("O7").Value = 0.3039081579394173
$ws.Range("P7").Value = 0.225443475773544
$ws.Range("Q7").Value = 1997.622456496291
$ws.Range("R7").Value = 11985.73473897775
$ws.Range("S7").Value = 0.04449331655521108
$ws.Range("T7").Value = 0.0340155735989895

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 123.444321
$ws.Range("H8").Value = 370.332963
$ws.Range("I8").Value = 0.2197321429647646
$ws.Range("J8").Value = 0.2264546783208506
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.275247666666667
$ws.Range("N8").Value = 24.825743
$ws.Range("O8").Value = 0.1035475654300768
$ws.Range("P8").Value = 0.1152196269592993
$ws.Range("Q8").Value = 1021.532329318501
$ws.Range("R8").Value = 9193.79096386651
$ws.Range("S8").Value = 0.02275272845073495
$ws.Range("T8").Value = 0.02609202355931652

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 123.444321
$ws.Range("H9").Value = 370.332963
$ws.Range("I9").Value = 0.2197321429647646
$ws.Range("J9").Value = 0.2264546783208506
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 47.354572
$ws.Range("N9").Value = 142.063716
$ws.Range("O9").Value = 0.5925442766305059
$ws.Range("P9").Value = 0.6593368972671567
$ws.Range("Q9").Value = 5845.652986785612
$ws.Range("R9").Value = 52610.87688107051
$ws.Range("S9").Value = 0.1302010237055273
$ws.Range("T9").Value = 0.1493099249757017

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 123.444321
$ws.Range("H10").Value = 370.332963
$ws.Range("I10").Value = 0.2197321429647646
$ws.Range("J10").Value = 0.2264546783208506
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 24.2875365
$ws.Range("N10").Value = 48.575073
$ws.Range("O10").Value = 0.3039081579394173
$ws.Range("P10").Value = 0.225443475773544
$ws.Range("Q10").Value = 2998.158452005217
$ws.Range("R10").Value = 17988.9507120313
$ws.Range("S10").Value = 0.0667783908085023
$ws.Range("T10").Value = 0.05105272978583238

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 154.3429766666667
$ws.Range("H11").Value = 463.02893
$ws.Range("I11").Value = 0.2747320633285943
$ws.Range("J11").Value = 0.2831372788071194
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.275247666666667
$ws.Range("N11").Value = 24.825743
$ws.Range("O11").Value = 0.1035475654300768
$ws.Range("P11").Value = 0.1152196269592993
$ws.Range("Q11").Value = 1277.226357527221
$ws.Range("R11").Value = 11495.03721774499
$ws.Range("S11").Value = 0.02844783630325762
$ws.Range("T11").Value = 0.03262297164242741

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 154.3429766666667
$ws.Range("H12").Value = 463.02893
$ws.Range("I12").Value = 0.2747320633285943
$ws.Range("J12").Value = 0.2831372788071194
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 47.354572
$ws.Range("N12").Value = 142.063716
$ws.Range("O12").Value = 0.5925442766305059
$ws.Range("P12").Value = 0.6593368972671567
$ws.Range("Q12").Value = 7308.845601255985
$ws.Range("R12").Value = 65779.61041130387
$ws.Range("S12").Value = 0.1627909117322482
$ws.Range("T12").Value = 0.186682854909352

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 154.3429766666667
$ws.Range("H13").Value = 463.02893
$ws.Range("I13").Value = 0.2747320633285943
$ws.Range("J13").Value = 0.2831372788071194
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 24.2875365
$ws.Range("N13").Value = 48.575073
$ws.Range("O13").Value = 0.3039081579394173
$ws.Range("P13").Value = 0.225443475773544
$ws.Range("Q13").Value = 3748.610679310314
$ws.Range("R13").Value = 22491.66407586189
$ws.Range("S13").Value = 0.08349331529308843
$ws.Range("T13").Value = 0.06383145225534

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 50.0323125
$ws.Range("H14").Value = 100.064625
$ws.Range("I14").Value = 0.08905802352064279
$ws.Range("J14").Value = 0.06118845668531954
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.275247666666667
$ws.Range("N14").Value = 24.825743
$ws.Range("O14").Value = 0.1035475654300768
$ws.Range("P14").Value = 0.1152196269592993
$ws.Range("Q14").Value = 414.0297772735626
$ws.Range("R14").Value = 2484.178663641375
$ws.Range("S14").Value = 0.009221741517577077
$ws.Range("T14").Value = 0.007050111153497759

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 50.0323125
$ws.Range("H15").Value = 100.064625
$ws.Range("I15").Value = 0.08905802352064279
$ws.Range("J15").Value = 0.06118845668531954
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 47.354572
$ws.Range("N15").Value = 142.063716
$ws.Range("O15").Value = 0.5925442766305059
$ws.Range("P15").Value = 0.6593368972671567
$ws.Range("Q15").Value = 2369.25874460775
$ws.Range("R15").Value = 14215.5524676465
$ws.Range("S15").Value = 0.05277082212518187
$ws.Range("T15").Value = 0.0403438071794644

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 50.0323125
$ws.Range("H16").Value = 100.064625
$ws.Range("I16").Value = 0.08905802352064279
$ws.Range("J16").Value = 0.06118845668531954
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 24.2875365
$ws.Range("N16").Value = 48.575073
$ws.Range("O16").Value = 0.3039081579394173
$ws.Range("P16").Value = 0.225443475773544
$ws.Range("Q16").Value = 1215.161616023156
$ws.Range("R16").Value = 4860.646464092625
$ws.Range("S16").Value = 0.02706545987788385
$ws.Range("T16").Value = 0.01379453835235738
